$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: "Save", matching the style of existing header cells (e.g. G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Data cells H2:H7: all zero, same (default) style as the other numeric data cells
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
